$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Rows 2-48: update Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = "63.290.01"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "3.255.58"
$ws.Range("E3").Value = "  +3.09%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "594.89"
$ws.Range("E5").Value = "  -1.54%  "
Set-TextValue "D6" "140.85"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.252.51"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("E10").Value = "  -1.37%  "
Set-TextValue "D11" "5.33"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "3.791.39"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "3.251.68"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").Value = "63.408.16"
$ws.Range("E18").Value = "  -1.14%  "
Set-TextValue "D19" "6.76"
$ws.Range("E19").Value = "  -1.39%  "
Set-TextValue "D20" "476.26"
$ws.Range("E20").Value = "  -3.01%  "
Set-TextValue "D21" "14.18"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("E22").Value = "  +1.77%  "
Set-TextValue "D23" "7.95"
$ws.Range("E23").Value = "  +3.79%  "
Set-TextValue "D24" "83.96"
$ws.Range("E24").Value = "  -4.57%  "
Set-TextValue "D25" "13.30"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -1.19%  "
Set-TextValue "D28" "7.43"
$ws.Range("E28").Value = "  +5.64%  "
Set-TextValue "D29" "8.08"
$ws.Range("E29").Value = "  -1.90%  "
Set-TextValue "D30" "2.12"
$ws.Range("E30").Value = "  +2.85%  "
Set-TextValue "D31" "27.67"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("E34").Value = "  -4.21%  "
Set-TextValue "D35" "1.09"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -2.42%  "
Set-TextValue "D37" "52.91"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "0.0₃0718"
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("E39").Value = "  -1.19%  "
Set-TextValue "D40" "420.48"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").Value = "3.000.91"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E43").Value = "  -7.28%  "
$ws.Range("E44").Value = "  -7.36%  "
Set-TextValue "D45" "0.269"
$ws.Range("E45").Value = "  +3.36%  "
Set-TextValue "D46" "2.16"
$ws.Range("E46").Value = "  -1.86%  "
Set-TextValue "D48" "25.90"
$ws.Range("E48").Value = "  -0.35%  "

# Rows 49-51: reorder coins (Stellar/ThetaToken/Monero -> ThetaToken/Stellar/Arweave) with updated data
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D49" "2.31"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.114"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D51" "33.76"
$ws.Range("E51").Value = "  +8.32%  "
